# Re-run k-means clustering results: refresh cluster id (col B) and the
# assigned specialty / therapeutic-labeling cluster names (cols C/D) for the
# rows whose cluster assignment changed, and append the two newly-clustered
# records (rows 85-86) that fell out of the updated k-means run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "other"
$ws.Range("D9").Value = "other"
$ws.Range("D26").Value = "cardiology"
$ws.Range("D32").Value = "other"
$ws.Range("D34").Value = "cardiology"
$ws.Range("B35").Value = 5
$ws.Range("D35").Value = "other"
$ws.Range("B36").Value = 5
$ws.Range("D37").Value = "cardiology"
$ws.Range("D38").Value = "other"
$ws.Range("D39").Value = "ocular"
$ws.Range("B40").Value = 6
$ws.Range("C40").Value = "nurse practioner"
$ws.Range("D40").Value = "bronchodilator"
$ws.Range("B41").Value = 6
$ws.Range("C41").Value = "nurse practioner"
$ws.Range("D41").Value = "other"
$ws.Range("D44").Value = "ocular"
$ws.Range("D45").Value = "ocular"
$ws.Range("B49").Value = 7
$ws.Range("D49").Value = "other"
$ws.Range("B50").Value = 7
$ws.Range("D50").Value = "other"
$ws.Range("D54").Value = "ocular"
$ws.Range("D55").Value = "ocular"
$ws.Range("B56").Value = 8
$ws.Range("C56").Value = "optometry"
$ws.Range("D56").Value = "other"
$ws.Range("B57").Value = 8
$ws.Range("C57").Value = "optometry"
$ws.Range("D57").Value = "other"
$ws.Range("D58").Value = "other"
$ws.Range("D59").Value = "other"
$ws.Range("B63").Value = 9
$ws.Range("B64").Value = 9
$ws.Range("D66").Value = "bronchodilator"
$ws.Range("B68").Value = 10
$ws.Range("C68").Value = "pulmonologist"
$ws.Range("D68").Value = "other"
$ws.Range("B69").Value = 10
$ws.Range("C69").Value = "pulmonologist"
$ws.Range("D69").Value = "bronchodilator"
$ws.Range("D71").Value = "ocular"
$ws.Range("D72").Value = "ocular"
$ws.Range("B73").Value = 11
$ws.Range("D73").Value = "other"
$ws.Range("B74").Value = 11
$ws.Range("D74").Value = "other"
$ws.Range("B78").Value = 12
$ws.Range("C78").Value = "unknown"
$ws.Range("D78").Value = "bronchodilator"
$ws.Range("B79").Value = 12
$ws.Range("C79").Value = "unknown"
$ws.Range("D79").Value = "bronchodilator"
$ws.Range("A80").Value = 79
$ws.Range("A81").Value = 80
$ws.Range("A82").Value = 81
$ws.Range("A83").Value = 82
$ws.Range("D83").Value = "cardiology"
$ws.Range("A84").Value = 83
$ws.Range("D84").Value = "cardiology"
$ws.Range("A85").Value = 84
$ws.Range("B85").Value = 13
$ws.Range("C85").Value = "cardiologist"
$ws.Range("D85").Value = "other"
$ws.Range("A86").Value = 85
$ws.Range("B86").Value = 13
$ws.Range("C86").Value = "cardiologist"
$ws.Range("D86").Value = "other"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F33").Select()

Write-Output "Applied clustering updates"
